$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 holds a date serial that needs to advance by one day.
$ws.Range("A1").Value = 45311

# Update the price column (D) for rows 29-37 with the new values.
$ws.Range("D29").Value = 62.043
$ws.Range("D30").Value = 66.429
$ws.Range("D31").Value = 69.563
$ws.Range("D32").Value = 73.01000000000001
$ws.Range("D33").Value = 76.459
$ws.Range("D34").Value = 80.84099999999999
$ws.Range("D35").Value = 87.738
$ws.Range("D36").Value = 99.01900000000001
$ws.Range("D37").Value = 115.94
